$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "63.473.13"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +2.22%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.554.20"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +5.37%  "
$ws.Range("E4").Value = "  +0.07%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "572.47"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.78%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "150.11"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +8.37%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.50%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.553.67"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +5.42%  "
$ws.Range("E10").Value = "  +2.15%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.76"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("E13").Value = "  +3.17%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "28.05"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +8.95%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "3.012.73"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +5.51%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "63.435.83"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.26%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.0000143"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +2.53%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.560.20"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +5.55%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.59"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +4.22%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "341.59"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("E21").Value = "  +2.92%  "
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("E23").Value = "  -0.01%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "66.15"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.74%  "
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("E26").Value = "  +3.75%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("E29").Value = "  +7.55%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "7.21"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +14.07%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.0₃0837"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +6.00%  "
$ws.Range("E32").Value = "  +4.03%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "177.53"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +3.96%  "
$ws.Range("E34").Value = "  +9.51%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "413.97"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +10.62%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.407"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +2.83%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "19.07"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.95%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "4.44"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("E39").Value = "  +0.01%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.75"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +4.53%  "
$ws.Range("E41").Value = "  -0.07%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "40.02"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.45%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "155.56"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +6.72%  "
$ws.Range("E44").Value = "  +3.39%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "21.04"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.77%  "
$ws.Range("E46").Value = "  +4.00%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0532"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.48%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0965"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("E49").Value = "  +5.66%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "18.73"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.83%  "
$ws.Range("E51").Value = "  +8.87%  "
